# Rename the inline picture shapes in the document's headers/footers.
#
# The Pearson Edexcel logo (alt text / descr:
#   "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png")
# appears in both the "first page" footer and the "default" footer and is
# renamed from image2.png -> image1.png.
#
# The BTEC logo (alt text / descr: "BTec_Logo-Orange") appears in the
# "first page" header and is renamed from image1.jpg -> image2.jpg.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($i = 1; $i -le 3; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        }
    }
}

Write-Output "Renamed header/footer logo images."
